# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the
# per-class Leve profit tables, as produced by the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 21737.5
$ws.Range("J3").Value = 21737.5
$ws.Range("L3").Value = 21737.5
$ws.Range("N3").Value = -21965.5

$ws.Range("H6").Value = 296
$ws.Range("I6").Value = 255.2
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 765.5999999999999
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -653.5999999999999
$ws.Range("N6").Value = -1724

$ws.Range("H28").Value = 845.35297
$ws.Range("I28").Value = 797.1429000000001
$ws.Range("K28").Value = 797.1429000000001
$ws.Range("M28").Value = -312.1429000000001

$ws.Range("H38").Value = 121.833336
$ws.Range("I38").Value = 42
$ws.Range("K38").Value = 126
$ws.Range("M38").Value = 246

$ws.Range("H39").Value = 194.6
$ws.Range("I39").Value = 102.22222
$ws.Range("J39").Value = 333.16666
$ws.Range("K39").Value = 306.66666
$ws.Range("L39").Value = 999.4999799999999
$ws.Range("M39").Value = -10.66665999999998
$ws.Range("N39").Value = -1591.49998

$ws.Range("H64").Value = 11221.667
$ws.Range("I64").Value = 8999.429
$ws.Range("K64").Value = 8999.429
$ws.Range("M64").Value = -8751.429

$ws.Range("H67").Value = 11221.667
$ws.Range("I67").Value = 8999.429
$ws.Range("K67").Value = 8999.429
$ws.Range("M67").Value = -8141.429

$ws.Range("H88").Value = 4047.875
$ws.Range("J88").Value = 4047.875
$ws.Range("L88").Value = 4047.875
$ws.Range("N88").Value = -4859.875

$ws.Range("H91").Value = 4047.875
$ws.Range("J91").Value = 4047.875
$ws.Range("L91").Value = 4047.875
$ws.Range("N91").Value = -6855.875

$ws.Range("H102").Value = 21737.5
$ws.Range("J102").Value = 21737.5
$ws.Range("L102").Value = 21737.5
$ws.Range("N102").Value = -28227.5

$ws.Range("H105").Value = 19000
$ws.Range("J105").Value = 19000
$ws.Range("L105").Value = 19000
$ws.Range("N105").Value = -25988

$ws.Range("H107").Value = 722.73334
$ws.Range("I107").Value = 614
$ws.Range("K107").Value = 614
$ws.Range("M107").Value = 1306

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 48533.3
$ws.Range("J24").Value = 48533.3
$ws.Range("L24").Value = 48533.3
$ws.Range("N24").Value = -49281.3

$ws.Range("H100").Value = 48533.3
$ws.Range("J100").Value = 48533.3
$ws.Range("L100").Value = 48533.3
$ws.Range("N100").Value = -50697.3

$ws.Range("H132").Value = 5310.857
$ws.Range("I132").Value = 3137.3333
$ws.Range("J132").Value = 6941
$ws.Range("K132").Value = 9411.999899999999
$ws.Range("L132").Value = 20823
$ws.Range("M132").Value = -6881.999899999999
$ws.Range("N132").Value = -25883

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2603.2
$ws.Range("J86").Value = 3835.6667
$ws.Range("L86").Value = 3835.6667
$ws.Range("N86").Value = -6081.6667

$ws.Range("H89").Value = 2603.2
$ws.Range("J89").Value = 3835.6667
$ws.Range("L89").Value = 19178.3335
$ws.Range("N89").Value = -30410.3335

$ws.Range("H100").Value = 9207.333000000001
$ws.Range("J100").Value = 9207.333000000001
$ws.Range("L100").Value = 9207.333000000001
$ws.Range("N100").Value = -11371.333

$ws.Range("H106").Value = 15000
$ws.Range("J106").Value = 15000
$ws.Range("L106").Value = 15000
$ws.Range("N106").Value = -17524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 15070.429
$ws.Range("J28").Value = 15070.429
$ws.Range("L28").Value = 15070.429
$ws.Range("N28").Value = -15560.429

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H99").Value = 4476.048
$ws.Range("I99").Value = 3945.2
$ws.Range("J99").Value = 5803.1665
$ws.Range("K99").Value = 3945.2
$ws.Range("L99").Value = 5803.1665
$ws.Range("M99").Value = -2447.2
$ws.Range("N99").Value = -8799.166499999999

$ws.Range("H122").Value = 749
$ws.Range("I122").Value = 499
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 1497
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = 953
$ws.Range("N122").Value = -7897

$ws.Range("H126").Value = 4476.048
$ws.Range("I126").Value = 3945.2
$ws.Range("J126").Value = 5803.1665
$ws.Range("K126").Value = 11835.6
$ws.Range("L126").Value = 17409.4995
$ws.Range("M126").Value = -9365.599999999999
$ws.Range("N126").Value = -22349.4995

$ws.Range("H132").Value = 5522.7
$ws.Range("I132").Value = 968.5454999999999
$ws.Range("K132").Value = 2905.6365
$ws.Range("M132").Value = -375.6364999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2923.6072
$ws.Range("J132").Value = 3499.65
$ws.Range("L132").Value = 31496.85
$ws.Range("N132").Value = -36556.85000000001

$ws.Range("H139").Value = 3449
$ws.Range("I139").Value = 1900
$ws.Range("J139").Value = 4998
$ws.Range("K139").Value = 5700
$ws.Range("L139").Value = 14994
$ws.Range("M139").Value = -560
$ws.Range("N139").Value = -25274

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1000000000
$ws.Range("I126").Value = 1000000000
$ws.Range("K126").Value = 3000000000
$ws.Range("M126").Value = -2999997530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8502
$ws.Range("I40").Value = 8502
$ws.Range("K40").Value = 8502
$ws.Range("M40").Value = -8366

$ws.Range("H61").Value = 1708.4667
$ws.Range("I61").Value = 1384.3334
$ws.Range("K61").Value = 1384.3334
$ws.Range("M61").Value = -1182.3334

$ws.Range("H82").Value = 2617.7144
$ws.Range("J82").Value = 2997.111
$ws.Range("L82").Value = 2997.111
$ws.Range("N82").Value = -3719.111

$ws.Range("H85").Value = 2617.7144
$ws.Range("J85").Value = 2997.111
$ws.Range("L85").Value = 2997.111
$ws.Range("N85").Value = -5493.111

$ws.Range("H113").Value = 1708.4667
$ws.Range("I113").Value = 1384.3334
$ws.Range("K113").Value = 1384.3334
$ws.Range("M113").Value = 785.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 34082.832
$ws.Range("J103").Value = 34082.832
$ws.Range("L103").Value = 34082.832
$ws.Range("N103").Value = -36426.832

$ws.Range("H107").Value = 620.6667
$ws.Range("I107").Value = 629.4286
$ws.Range("J107").Value = 613
$ws.Range("K107").Value = 1888.2858
$ws.Range("L107").Value = 1839
$ws.Range("M107").Value = 31.71420000000012
$ws.Range("N107").Value = -5679

$ws.Range("H113").Value = 547.4
$ws.Range("I113").Value = 387.75
$ws.Range("J113").Value = 653.8333
$ws.Range("K113").Value = 1163.25
$ws.Range("L113").Value = 1961.4999
$ws.Range("M113").Value = 1006.75
$ws.Range("N113").Value = -6301.4999

$ws.Range("H126").Value = 3779
$ws.Range("I126").Value = 3473.75
$ws.Range("K126").Value = 10421.25
$ws.Range("M126").Value = -7951.25

$ws.Range("H141").Value = 95333
$ws.Range("J141").Value = 95333
$ws.Range("L141").Value = 95333
$ws.Range("N141").Value = -105693
